{"js": "// The author rewrote the opening of the abstract's only paragraph:\n//   1) \"In order to design the spacesuit vitals portion of the project we\n//       will create a program that\" -> \"The spacesuit vitals application\"\n//   2) \"a call is immediately initiated with ground control.\" ->\n//       \"contact with ground control is immediately established.\"\n// Everything else in the paragraph (including the \"_GoBack\" bookmark and\n// the spell-check proofErr markers around \"vitals\") is left untouched.\n\nconst body = context.document.body;\n\n// --- Edit 1 --------------------------------------------------------------\nlet results = body.search(\n  \"In order to design the spacesuit vitals portion of the project we will create a program that\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"The spacesuit vitals application\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2 --------------------------------------------------------------\nresults = body.search(\n  \"a call is immediately initiated with ground control.\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"contact with ground control is immediately established.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# The author rewrote the opening of the abstract's only paragraph:\n#   1) \"In order to design the spacesuit vitals portion of the project we\n#       will create a program that\" -> \"The spacesuit vitals application\"\n#   2) \"a call is immediately initiated with ground control.\" ->\n#       \"contact with ground control is immediately established.\"\n# Everything else in the paragraph (including the \"_GoBack\" bookmark and\n# the spell-check proofErr markers around \"vitals\") is left untouched.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Replacement.ClearFormatting()\n$range1.Find.Execute(\n    \"In order to design the spacesuit vitals portion of the project we will create a program that\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"The spacesuit vitals application\",\n    2\n) | Out-Null\n\n# --- Edit 2 ---------------------------------------------------------------\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Replacement.ClearFormatting()\n$range2.Find.Execute(\n    \"a call is immediately initiated with ground control.\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"contact with ground control is immediately established.\",\n    2\n) | Out-Null\n"}
